{"js": "// Replace the \"incident of violence, threat of violence, or other abusive\n// behavior by\" phrasing with \"incident of stalking or other threatening\n// behavior by\" in the statement-of-facts template. Per the authored change,\n// the sentence fragment is rebuilt as three separate runs of text (all\n// sharing the original run's formatting: black color, 10pt):\n//   1) \"%} incident of \"\n//   2) \"stalking or other threatening behavior by\"\n//   3) \" {{\"\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\n  \"%} incident of violence, threat of violence, or other abusive behavior by {{\",\n  { matchCase: true }\n);\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Target phrase not found in document.\");\n}\n\nconst target = searchResults.items[0];\n\n// Capture the original run's formatting before it gets overwritten, so the\n// newly inserted runs can be given the same look explicitly.\ntarget.font.load(\"color,size,highlightColorIndex\");\nawait context.sync();\n\nconst originalColor = target.font.color;\nconst originalSize = target.font.size;\n\nfunction applyFormatting(range) {\n  range.font.color = originalColor;\n  range.font.size = originalSize;\n}\n\n// First chunk replaces the whole matched range ...\nconst run1 = target.insertText(\"%} incident of \", Word.InsertLocation.replace);\napplyFormatting(run1);\nawait context.sync();\n\n// ... then the remaining two chunks are appended right after it, in order.\nconst run2 = run1.insertText(\"stalking or other threatening behavior by\", Word.InsertLocation.after);\napplyFormatting(run2);\nawait context.sync();\n\nconst run3 = run2.insertText(\" {{\", Word.InsertLocation.after);\napplyFormatting(run3);\nawait context.sync();\n", "ps1": "# Replace the \"incident of violence, threat of violence, or other abusive\n# behavior by\" phrasing with \"incident of stalking or other threatening\n# behavior by\" in the statement-of-facts template. Per the authored change,\n# the sentence fragment is rebuilt as three separate runs of text (all\n# sharing the original run's formatting: black color, 10pt):\n#   1) \"%} incident of \"\n#   2) \"stalking or other threatening behavior by\"\n#   3) \" {{\"\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"%} incident of violence, threat of violence, or other abusive behavior by {{\"\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Target phrase not found in document.\"\n}\n\n$matchStart = $find.Parent.Start\n$matchEnd = $find.Parent.End\n\n$target = $d.Range($matchStart, $matchEnd)\n\n# Remember the matched run's formatting so the rebuilt runs keep it.\n$origColor = $target.Font.Color\n$origSize = $target.Font.Size\n\n$part1 = \"%} incident of \"\n$part2 = \"stalking or other threatening behavior by\"\n$part3 = \" {{\"\n\n# Step 1: collapse the whole matched range down to the first chunk of text.\n$target.Text = $part1\n$target.Font.Color = $origColor\n$target.Font.Size = $origSize\n\n# Step 2: insert the second chunk right after the first, as its own run.\n$insertPoint2 = $matchStart + $part1.Length\n$range2 = $d.Range($insertPoint2, $insertPoint2)\n$range2.InsertAfter($part2)\n$range2 = $d.Range($insertPoint2, $insertPoint2 + $part2.Length)\n$range2.Font.Color = $origColor\n$range2.Font.Size = $origSize\n\n# Step 3: insert the third chunk right after the second, as its own run.\n$insertPoint3 = $insertPoint2 + $part2.Length\n$range3 = $d.Range($insertPoint3, $insertPoint3)\n$range3.InsertAfter($part3)\n$range3 = $d.Range($insertPoint3, $insertPoint3 + $part3.Length)\n$range3.Font.Color = $origColor\n$range3.Font.Size = $origSize\n\nWrite-Output \"done\"\n"}
